# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46006
$ws.Range("B2").Value = 83.43000000000001
$ws.Range("C2").Value = 71.73
$ws.Range("D2").Value = 71.18000000000001
$ws.Range("E2").Value = 70.12
$ws.Range("F2").Value = 67.03
$ws.Range("G2").Value = 70.98
$ws.Range("H2").Value = 83.11
$ws.Range("I2").Value = 96.58
$ws.Range("J2").Value = 106.1
$ws.Range("K2").Value = 101.33
$ws.Range("L2").Value = 95.59999999999999
$ws.Range("M2").Value = 88.56
$ws.Range("N2").Value = 84.65000000000001
$ws.Range("O2").Value = 85.44
$ws.Range("P2").Value = 91.2
$ws.Range("Q2").Value = 99.39
$ws.Range("R2").Value = 104.35
$ws.Range("S2").Value = 104.26
$ws.Range("T2").Value = 102.52
$ws.Range("U2").Value = 105.43
$ws.Range("V2").Value = 115.04
$ws.Range("W2").Value = 112.54
$ws.Range("X2").Value = 109.77
$ws.Range("Y2").Value = 98.97
$ws.Range("Z2").Value = 92.47
$ws.Range("AB2").Value = 109.08
$ws.Range("AD2").Value = 113.79
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 104.37
$ws.Range("AG2").Value = "0h-14h"
